$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Row 7 = "Experimental" property. It previously had no Value; the FHIR
# IG publisher now always emits this element, so fill in the literal text
# "true" (not the Excel boolean TRUE) in the previously-empty cell.
# Writing the literal word via .Value/.Formula gets auto-coerced to a
# logical (boolean) cell by Excel's type inference, so instead compute it
# as a formula returning the text "true" and then collapse it down to a
# plain text value with Paste Special (Values) - this keeps the cell's
# original style/format intact while still storing real text.
$cell = $ws.Range("B7")
$cell.Formula = "=""true"""
$cell.Copy() | Out-Null
$cell.PasteSpecial(-4163) | Out-Null   # xlPasteValues
$excel.CutCopyMode = $false

# Row 8 = "Date" property; bump the timestamp to the new commit date.
$ws.Range("B8").Value = "2023-02-01T09:05:11-06:00"
